$wb = $excel.ActiveWorkbook

# "Planilha1" is the template sheet for the account/column description
# table. Duplicate it (Excel names the copy "Planilha1 (2)" automatically
# and places it right after the source sheet, becoming the new active
# sheet/tab) to create the new dim_contas description sheet.
$source = $wb.Worksheets.Item("Planilha1")
$source.Copy([System.Reflection.Missing]::Value, $source) | Out-Null

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count())

# Update the table title in B2 from "Tabela: PlanoContas" to the new
# "Tabela: dim_contas" label - everything else (headers, rows, styles)
# stays identical to the Planilha1 template.
$newSheet.Range("B2").Value = "Tabela: dim_contas"

# Match the author's final selection on the new sheet.
$newSheet.Range("B12").Select() | Out-Null
